# Adds "Nour Aldin Almubarak" to the Collaborative Filtering credit lines
# on slide 3 and slide 7 of the Content Placeholder text boxes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3: "Collaborative Filtering (Library) – Victor Essien"
#   becomes: "Collaborative Filtering – Victor Essien, Nour Aldin Almubarak"
#   (the existing run's text changes, then the run is split so the new
#   "Aldin" / " " / "Almubarak" pieces live in their own runs, matching
#   the authoring pattern used elsewhere in the deck)
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(5, 1)

# Replace the old "(Library) – Victor Essien" run's text with the new,
# longer text for the whole tail of the line.
$run3 = $para3.Characters(25, 25)
$run3.Text = "– Victor Essien, Nour Aldin Almubarak"

# Re-assert formatting on each logical piece so the single grown run is
# split back into separate runs at the correct boundaries.
$part_dash = $para3.Characters(25, 22)    # "– Victor Essien, Nour "
$part_dash.Font.Bold = $true
$part_dash.Font.Italic = $true

$part_aldin = $para3.Characters(47, 5)    # "Aldin"
$part_aldin.Font.Bold = $true
$part_aldin.Font.Italic = $true

$part_space = $para3.Characters(52, 1)    # " "
$part_space.Font.Bold = $true
$part_space.Font.Italic = $true

$part_almubarak = $para3.Characters(53, 9) # "Almubarak"
$part_almubarak.Font.Bold = $true
$part_almubarak.Font.Italic = $true

# ---------------------------------------------------------------------------
# Slide 7: "Collaborative Filtering (Library) – Victor Essien, Nour Aldin"
#   becomes: "...Nour Aldin Almubarak"
#   Two new runs ("Aldin" and " ") are inserted before the old "Aldin" run,
#   which itself becomes "Almubarak". The new runs are carved out of the
#   "(Library) – Victor Essien, Nour " run (grown to include "Aldin ") so
#   they naturally pick up that run's existing formatting.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$tr7 = $shp7.TextFrame.TextRange
$para7 = $tr7.Paragraphs(1, 1)

# Grow the "(Library) – Victor Essien, Nour " run to also contain "Aldin ".
$run7_prefix = $para7.Characters(25, 32)
$run7_prefix.Text = "(Library) – Victor Essien, Nour Aldin "

# The old "Aldin" run has shifted right by 6 characters ("Aldin " minus
# nothing removed); rename its text to "Almubarak".
$run7_tail = $para7.Characters(63, 5)
$run7_tail.Text = "Almubarak"

# Re-assert formatting on the three logical pieces so the grown run is
# split back into "(Library) – Victor Essien, Nour " / "Aldin" / " ".
$part_prefix7 = $para7.Characters(25, 32)   # "(Library) – Victor Essien, Nour "
$part_prefix7.Font.Bold = $true
$part_prefix7.Font.Italic = $true

$part_aldin7 = $para7.Characters(57, 5)     # "Aldin"
$part_aldin7.Font.Bold = $true
$part_aldin7.Font.Italic = $true

$part_space7 = $para7.Characters(62, 1)     # " "
$part_space7.Font.Bold = $true
$part_space7.Font.Italic = $true
